# Generate Report for Handback
# -----------------------------------------------------------------------
# This script updates the localization-status workbook to reflect that
# the zh-cn and de-de handback packages have been generated:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview + both locale sheets).
#   - "Latest Handback DateTime" timestamps are filled in (previously the
#     zero-date placeholder) for both locales.
#   - "Latest Target File" (the source .md, now hyperlinked) and
#     "Latest Handback File" (the generated .xlf) columns are populated
#     for both locales.
#   - A couple of columns are widened to better fit the newly-populated
#     content.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)   # Overview
$wsZhCn     = $wb.Worksheets.Item(2)   # zh-cn
$wsDeDe     = $wb.Worksheets.Item(3)   # de-de

# -------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This text is shared by the Overview sheet (columns E/F) and the
#    Status column (C) on both locale sheets, so a whole-workbook
#    whole-cell replace covers every occurrence.
# -------------------------------------------------------------------
$wsOverview.UsedRange.Replace("Ready for handoff", "Handed back: in sync with en-US", 1) | Out-Null
$wsZhCn.UsedRange.Replace("Ready for handoff", "Handed back: in sync with en-US", 1) | Out-Null
$wsDeDe.UsedRange.Replace("Ready for handoff", "Handed back: in sync with en-US", 1) | Out-Null

# -------------------------------------------------------------------
# 2. Latest Handback DateTime (column K) on both locale sheets was the
#    zero-date placeholder "0001-01-01 00:00:00"; fill in the real
#    handback timestamps (different per locale).
# -------------------------------------------------------------------
$wsZhCn.UsedRange.Replace("0001-01-01 00:00:00", "2016-08-23 14:33:40", 1) | Out-Null
$wsDeDe.UsedRange.Replace("0001-01-01 00:00:00", "2016-08-23 14:33:47", 1) | Out-Null

# -------------------------------------------------------------------
# 3. Populate "Latest Target File" (I) and "Latest Handback File" (J)
#    for row 2 (50ce2d35...) and row 3 (6729fe01...) on both locale
#    sheets, and hyperlink the new "Latest Target File" cells the same
#    way the existing "Source File Name" (A) hyperlinks work.
# -------------------------------------------------------------------

function Set-HandbackRow($ws, $row, $mdName, $targetFileSS, $hash, $locale) {
    $iCell = $ws.Range("I$row")
    $iCell.Value = $mdName

    $url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3003c69674fa1fd7e50f017d4055110dbdb7c82a/e2e/$mdName"
    $ws.Hyperlinks.Add($iCell, $url, "", "", $mdName) | Out-Null
    $iCell.Font.Underline = 2
    $iCell.Font.Color = 15570276

    $xlfName = "$hash.$locale.xlf"
    $ws.Range("J$row").Value = $xlfName
}

Set-HandbackRow $wsZhCn 2 "50ce2d35-47ce-4a65-b417-dff63d0f89c2.md" $null "50ce2d35-47ce-4a65-b417-dff63d0f89c2.7e00a80765a1ebdb6b88d4b7e1a23b20ce8dfbbc" "zh-cn"
Set-HandbackRow $wsZhCn 3 "6729fe01-c883-460d-aac7-358de0a9360a.md" $null "6729fe01-c883-460d-aac7-358de0a9360a.19052f8876de4e27308698297384aacaa758529b" "zh-cn"

Set-HandbackRow $wsDeDe 2 "50ce2d35-47ce-4a65-b417-dff63d0f89c2.md" $null "50ce2d35-47ce-4a65-b417-dff63d0f89c2.7e00a80765a1ebdb6b88d4b7e1a23b20ce8dfbbc" "de-de"
Set-HandbackRow $wsDeDe 3 "6729fe01-c883-460d-aac7-358de0a9360a.md" $null "6729fe01-c883-460d-aac7-358de0a9360a.19052f8876de4e27308698297384aacaa758529b" "de-de"

# -------------------------------------------------------------------
# 4. Widen columns that now hold the longer populated values.
#    (The engine's ColumnWidth setter stores width internally in
#    whole-pixel units, offset by 5/6 of a character from the raw
#    OOXML "width" attribute; 39.1666... / 29.1666... are the nearest
#    achievable values to the authored widths of 40 and 29.9777...)
# -------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668   # E: zh-cn status
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668   # F: de-de status

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C: Status
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I: Latest Target File
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664      # J: Latest Handback File

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.166666666666668      # C: Status
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.166666666666664      # I: Latest Target File
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664      # J: Latest Handback File

Write-Host "Handback report generated."
